# commiting directly with latest changes where generic xpaths has been added
# to Dynamic Rule and INDB Creation

$wb = $excel.ActiveWorkbook

# --- Sheet1: INDBID -----------------------------------------------------
# Append newly discovered INDB ids below the existing row.
$wsIndbId = $wb.Worksheets.Item("INDBID")
$wsIndbId.Range("A3").Value = "onlineFPtjCXIAH2"
$wsIndbId.Range("A4").Value = "onlinecV6KDHIAAK"
$wsIndbId.Range("A5").Value = "onlineKDsdDnIAuF"

# --- Sheet3: REV_INDBID --------------------------------------------------
# Add the matching revision row (old/new INDB + null INDB ids).
$wsRevIndbId = $wb.Worksheets.Item("REV_INDBID")
$wsRevIndbId.Range("F3").Value = "onlineddiZynEAKQ"
$wsRevIndbId.Range("G3").Value = "onlineUVKcynEA1Q"
$wsRevIndbId.Range("H3").Value = "onlineZQ2aynEAHI"

# --- Sheet4: COPY_INDB -----------------------------------------------------
# Move the stored selection; this sheet is no longer the active tab.
$wsCopyIndb = $wb.Worksheets.Item("COPY_INDB")
$wsCopyIndb.Range("D14").Select()

# --- Sheet2: NullINDB ------------------------------------------------------
# Replace the stale null-INDB id and append the newly generated ones, then
# leave this sheet active/selected (matches activeTab moving to index 1).
$wsNullIndb = $wb.Worksheets.Item("NullINDB")
$wsNullIndb.Range("A2").Value = "onlineAS1DynEAFw"
$wsNullIndb.Range("A3").Value = "onlineWOxEynEAGw"
$wsNullIndb.Range("A4").Value = "onlineddiZynEAKQ"
$wsNullIndb.Select()
$wsNullIndb.Range("C8").Select()
